$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the used range extent
$used = $ws.UsedRange
$lastRow = $used.Rows.Count()

# Columns C (3) and D (4) - "codeforiati:group-name" and "codeforiati:group-code" -
# need to be swapped for the header and every data row, so that the
# group-code column comes before the group-name column.
for ($r = 1; $r -le $lastRow; $r++) {
    $cVal = $ws.Cells.Item($r, 3).Value()
    $dVal = $ws.Cells.Item($r, 4).Value()
    $ws.Cells.Item($r, 3).Value = $dVal
    $ws.Cells.Item($r, 4).Value = $cVal
}
